$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Ray Tracing" overview): keep 1st bullet, change the other three
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "Essential Ray Tracing Algorithms`rRay Tracing Shaders`rRay Intersection Algorithms`rRay Object Traversal"

# ---------------------------------------------------------------------------
# Slide 3 ("Overview" -> "Ray Tracing Shaders"); bullets -> shader list
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Text = "Ray"
[void]$title3.InsertAfter(" ")
[void]$title3.InsertAfter("Tracing")
[void]$title3.InsertAfter(" ")
[void]$title3.InsertAfter("Shaders")

$body3 = $s3.Shapes.Item(2).TextFrame.TextRange
$body3.Text = "Ray Generation Shader`rMiss Shader`rClosest Hit Shader`rAny Hit Shader`rBounding Volume Hierarchy Processor"

# ---------------------------------------------------------------------------
# Slide 4 ("Ray Tracing Shaders" -> "Ray Intersection Algorithms");
# content -> Ray/Plane intersection bullets
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$title4 = $s4.Shapes.Item(1).TextFrame.TextRange
$title4.Text = "Ray"
[void]$title4.InsertAfter(" ")
[void]$title4.InsertAfter("Intersection")
[void]$title4.InsertAfter(" ")
[void]$title4.InsertAfter("Algorithms")

$body4 = $s4.Shapes.Item(2).TextFrame.TextRange
$body4.Text = "Ray Sphere Intersection`rRay Plane Intersection"

# ---------------------------------------------------------------------------
# New slide (appended at the end) - "Ray Object Traversal" with pseudocode
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Add($p.Slides.Count + 1, 2)
$title5 = $s5.Shapes.Item(1).TextFrame.TextRange
$title5.Text = "Ray"
[void]$title5.InsertAfter(" ")
[void]$title5.InsertAfter("Object")
[void]$title5.InsertAfter(" ")
[void]$title5.InsertAfter("Traversal")

$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = "h <- HITRECORD(cur_t = tmax)`rLoop through object list`r    t <- INTERSECT(ray, object)`r    IF t < tmin THEN CONTINUE`r    IF t > tmax THEN CONTINUE`r    IF t > cur_t THEN CONTINUE`r    cur_t <- t`r    ANYHIT_SHADER(h)`rEND LOOP`rIF cur_t != tmax THEN CLOSEST_HIT_SHADER(h)`rELSE MISS_SHADER(h)"
$body5.Font.Size = 18
$body5.Font.Name = "Courier"
